# UndoAndRedoImplementation.pptx edit script
#
# 1) Refresh the cached "datetimeFigureOut" field text on the slide master
#    and all 11 slide layouts from "1/11/2019" to "9/11/2019".
# 2) On slide 5, move the "a:AddExpenseCommand" rectangle (Rectangle 9) down
#    to y=5533458 EMU and remove the now-redundant "u:UpdateCommand"
#    rectangle (Rectangle 14).
# 3) On slide 6, remove the same redundant "u:UpdateCommand" rectangle
#    (Rectangle 14).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text refresh (master + every layout)
# ---------------------------------------------------------------------------
$newDate = "9/11/2019"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 5: reposition Rectangle 9, delete Rectangle 14
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)

for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $sh = $slide5.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 9") {
        # 5533458 EMU, expressed in points (1 pt = 12700 EMU)
        $sh.Top = 435.70538330078125
    }
}

for ($i = $slide5.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide5.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 14") {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 6: delete Rectangle 14
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

for ($i = $slide6.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide6.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 14") {
        $sh.Delete()
    }
}
